# SwaadSutra_Daily_2026-01-21.xlsx update
# Adds one new order row to "Daily Orders", rolls the totals into
# "Summary", and rebuilds the per-item "Items Breakdown".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Daily Orders
# ---------------------------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Daily Orders")

$headers = @("Order ID","Date","Customer","Flat No","Phone","Items","Total","Status","Payment","Collection Date","Collection Time","Notes","Cancel Reason","Feedback")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsOrders.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Leading apostrophe = Excel's own "force text / quote-prefix" convention,
# needed so numeric/date-shaped text (phone number, ISO date) is kept as a
# literal string instead of being auto-coerced into a Number/Date.
# ClearFormats() afterwards drops the quote-prefix style flag again (the
# value stays text) so no stray cell style is left behind.
$wsOrders.Range("A2").Value = 27
$wsOrders.Range("B2").Value = "2026-01-21 07:49"
$wsOrders.Range("C2").Value = "Renu"
$wsOrders.Range("D2").Value = "A-1005 Kakkad la vida"
$wsOrders.Range("E2").Value = "'8806022013"
$wsOrders.Range("E2").ClearFormats()
$wsOrders.Range("F2").Value = "Appe Chutney x1, Vermicelli Kheer x1"
$wsOrders.Range("G2").Value = 110
$wsOrders.Range("H2").Value = "NEW"
$wsOrders.Range("I2").Value = "PENDING"
$wsOrders.Range("J2").Value = "'2026-01-21"
$wsOrders.Range("J2").ClearFormats()
$wsOrders.Range("K2").Value = "18:30"
$wsOrders.Range("L2").Value = "Less spicy"
# A lone quote-prefix character creates a *static* empty-text cell (as
# opposed to assigning "" outright, which Excel/this engine treats as
# "clear the cell" and would leave no cell behind at all).
$wsOrders.Range("M2").Value = "'"
$wsOrders.Range("M2").ClearFormats()
$wsOrders.Range("N2").Value = "'"
$wsOrders.Range("N2").ClearFormats()

# ---------------------------------------------------------------------
# Sheet 2: Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("A2").Value = 1    # Total Orders
$wsSummary.Range("B2").Value = 1    # New
$wsSummary.Range("C2").Value = 0    # Cooking
$wsSummary.Range("D2").Value = 0    # Ready
$wsSummary.Range("E2").Value = 0    # Delivered
$wsSummary.Range("F2").Value = 0    # Cancelled
$wsSummary.Range("G2").Value = 110  # Total Revenue
$wsSummary.Range("H2").Value = 0    # Paid Amount

# ---------------------------------------------------------------------
# Sheet 3: Items Breakdown
# ---------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("Items Breakdown")

$wsItems.Range("A1").Value = "Item"
$wsItems.Range("B1").Value = "Quantity Ordered"
$wsItems.Range("C1").Value = "Revenue"

$wsItems.Range("A2").Value = "Appe Chutney"
$wsItems.Range("B2").Value = 1
$wsItems.Range("C2").Value = 60

$wsItems.Range("A3").Value = "Vermicelli Kheer"
$wsItems.Range("B3").Value = 1
$wsItems.Range("C3").Value = 50
